$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: weight halved (approx), add a Numb value
$ws.Range("G3").Value = 4.2
$ws.Range("H3").Value = 492

# Row 23: weight halved (approx), add a Numb value
$ws.Range("G23").Value = 4.062
$ws.Range("H23").Value = 426

# Row 27: small weight correction
$ws.Range("G27").Value = 1.813

# Rows 38-61: RF (column I) recalculated to a new raising-factor value
$newRF = 99.11028571428572
for ($r = 38; $r -le 61; $r++) {
    $ws.Cells.Item($r, 9).Value = $newRF
}

# Rows with Numb (column H) corrected from 0 to -1
$ws.Range("H46").Value = -1
$ws.Range("H60").Value = -1
$ws.Range("H61").Value = -1
